$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data. Columns D (Price) values that look numeric
# must be forced to Text format so Excel keeps them as literal strings
# (matching the source inline-string cells, e.g. "89.300.77", "0.0700").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '89.300.77'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.091.32'
$ws.Range("E3").Value = '  -2.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.78'
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '622.07'
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.371'
$ws.Range("E7").Value = '  -6.05%  '
$ws.Range("E8").Value = '  +15.58%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.087.32'
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.614'
$ws.Range("E11").Value = '  +9.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.182'
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("E13").Value = '  -4.40%  '
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.089.14'
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.19'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.663.29'
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.094.68'
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("E19").Value = '  +3.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000213'
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("E22").Value = '  -2.97%  '
$ws.Range("E23").Value = '  -2.22%  '
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.52'
$ws.Range("E25").Value = '  +6.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.09'
$ws.Range("E26").Value = '  +4.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '83.58'
$ws.Range("E27").Value = '  +4.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.252.95'
$ws.Range("E28").Value = '  -2.53%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("B30").Value = 'Cronos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.165'
$ws.Range("E30").Value = '  +6.55%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.01'
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '508.56'
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("E34").Value = '  -7.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.69'
$ws.Range("E35").Value = '  -3.21%  '
$ws.Range("E36").Value = '  -1.85%  '
$ws.Range("E37").Value = '  -3.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.40'
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -1.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.136'
$ws.Range("E44").Value = '  +8.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.83'
$ws.Range("E45").Value = '  -4.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '146.34'
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0700'
$ws.Range("E47").Value = '  +13.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.51'
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("E49").Value = '  +1.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '159.30'
$ws.Range("E50").Value = '  -6.74%  '
$ws.Range("E51").Value = '  -4.66%  '
